$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99. This shifts the existing rows 99-119
# down to 100-120, preserving their data/formatting (mirrors a new weekly
# price record being inserted chronologically before the previous entries).
$ws.Rows.Item(99).Insert()

# The new blank row 99 should start as a copy of the row beneath it (which
# now holds what used to be row 99), then have its price-specific fields
# updated to the new observation.
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item(99, $col).Value = $ws.Cells.Item(100, $col).Value()
}

# Apply the new values for the freshly inserted record (row 99)
$ws.Cells.Item(99, 4).Value = 45173    # D99 - Fecha
$ws.Cells.Item(99, 13).Value = 150     # M99 - Volumen
$ws.Cells.Item(99, 14).Value = 24000   # N99 - Precio mínimo
$ws.Cells.Item(99, 15).Value = 24000   # O99 - Precio máximo
$ws.Cells.Item(99, 16).Value = 24000   # P99 - Precio promedio ponderado
$ws.Cells.Item(99, 19).Value = 2400    # S99 - Precio $/Kg
